$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.215933680534363
$ws.Range("B1").Value = 2.518023014068604
$ws.Range("D1").Value = 1.407598376274109
$ws.Range("E1").Value = 0.9150463342666626
